# Auto-generated script applying cell value updates per the commit diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1050.7931
$ws.Range("I135").Value = 748.43475
$ws.Range("J135").Value = 2209.8333
$ws.Range("K135").Value = 6735.91275
$ws.Range("L135").Value = 19888.4997
$ws.Range("M135").Value = -4200.91275
$ws.Range("N135").Value = -24958.4997
$ws.Range("H138").Value = 7357.3696
$ws.Range("I138").Value = 1526.3077
$ws.Range("J138").Value = 14937.75
$ws.Range("K138").Value = 4578.9231
$ws.Range("L138").Value = 44813.25
$ws.Range("M138").Value = 561.0769
$ws.Range("N138").Value = -55093.25
$ws.Range("H141").Value = 1924.4242
$ws.Range("I141").Value = 1686.2413
$ws.Range("J141").Value = 3651.25
$ws.Range("K141").Value = 5058.7239
$ws.Range("L141").Value = 10953.75
$ws.Range("M141").Value = 121.2761
$ws.Range("N141").Value = -21313.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2506.1614
$ws.Range("I61").Value = 1869.3077
$ws.Range("J61").Value = 2966.111
$ws.Range("K61").Value = 1869.3077
$ws.Range("L61").Value = 2966.111
$ws.Range("M61").Value = -1657.3077
$ws.Range("N61").Value = -3390.111
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("H74").Value = 2173.4666
$ws.Range("I74").Value = 852.875
$ws.Range("J74").Value = 3682.7144
$ws.Range("K74").Value = 852.875
$ws.Range("L74").Value = 3682.7144
$ws.Range("M74").Value = 21.125
$ws.Range("N74").Value = -5430.7144
$ws.Range("H77").Value = 2173.4666
$ws.Range("I77").Value = 852.875
$ws.Range("J77").Value = 3682.7144
$ws.Range("K77").Value = 4264.375
$ws.Range("L77").Value = 18413.572
$ws.Range("M77").Value = 103.625
$ws.Range("N77").Value = -27149.572
$ws.Range("H102").Value = 52100.5
$ws.Range("I102").Value = 72789.28999999999
$ws.Range("K102").Value = 72789.28999999999
$ws.Range("M102").Value = -71167.28999999999
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("H132").Value = 2390.827
$ws.Range("I132").Value = 2218.3333
$ws.Range("J132").Value = 3499.7144
$ws.Range("K132").Value = 6654.999899999999
$ws.Range("L132").Value = 10499.1432
$ws.Range("M132").Value = -4124.999899999999
$ws.Range("N132").Value = -15559.1432
$ws.Range("H136").Value = 2506.1614
$ws.Range("I136").Value = 1869.3077
$ws.Range("J136").Value = 2966.111
$ws.Range("K136").Value = 5607.9231
$ws.Range("L136").Value = 8898.332999999999
$ws.Range("M136").Value = -3057.9231
$ws.Range("N136").Value = -13998.333
$ws.Range("N70").ClearContents()
$ws.Range("N73").ClearContents()
$ws.Range("N125").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2811.111
$ws.Range("I134").Value = 2900.8484
$ws.Range("J134").Value = 2564.3333
$ws.Range("K134").Value = 8702.5452
$ws.Range("L134").Value = 7692.999899999999
$ws.Range("M134").Value = -6167.5452
$ws.Range("N134").Value = -12762.9999
$ws.Range("H135").Value = 46390
$ws.Range("J135").Value = 46390
$ws.Range("L135").Value = 46390
$ws.Range("N135").Value = -56530
$ws.Range("H137").Value = 40780
$ws.Range("J137").Value = 40780
$ws.Range("L137").Value = 40780
$ws.Range("N137").Value = -50980
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("N139").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5312.6484
$ws.Range("I58").Value = 1228.9667
$ws.Range("K58").Value = 1228.9667
$ws.Range("M58").Value = -1025.9667
$ws.Range("H70").Value = 14333.333
$ws.Range("J70").Value = 14333.333
$ws.Range("L70").Value = 14333.333
$ws.Range("N70").Value = -14963.333
$ws.Range("H73").Value = 14333.333
$ws.Range("J73").Value = 14333.333
$ws.Range("L73").Value = 14333.333
$ws.Range("N73").Value = -16517.333
$ws.Range("H105").Value = 2425.4167
$ws.Range("I105").Value = 2515.158
$ws.Range("J105").Value = 2084.4
$ws.Range("K105").Value = 2515.158
$ws.Range("L105").Value = 2084.4
$ws.Range("M105").Value = -768.1579999999999
$ws.Range("N105").Value = -5578.4
$ws.Range("H134").Value = 1366.4
$ws.Range("I134").Value = 1326.5
$ws.Range("K134").Value = 3979.5
$ws.Range("M134").Value = -1444.5
$ws.Range("H136").Value = 5312.6484
$ws.Range("I136").Value = 1228.9667
$ws.Range("K136").Value = 3686.9001
$ws.Range("M136").Value = -1136.9001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 4203.375
$ws.Range("I109").Value = 2900
$ws.Range("J109").Value = 4389.5713
$ws.Range("K109").Value = 8700
$ws.Range("L109").Value = 13168.7139
$ws.Range("M109").Value = -7660
$ws.Range("N109").Value = -15248.7139
$ws.Range("H113").Value = 632.80646
$ws.Range("I113").Value = 553
$ws.Range("J113").Value = 676.7
$ws.Range("K113").Value = 1659
$ws.Range("L113").Value = 2030.1
$ws.Range("M113").Value = 511
$ws.Range("N113").Value = -6370.1

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3127.3044
$ws.Range("I132").Value = 1895.3846
$ws.Range("J132").Value = 4728.8
$ws.Range("K132").Value = 5686.1538
$ws.Range("L132").Value = 14186.4
$ws.Range("M132").Value = -3156.1538
$ws.Range("N132").Value = -19246.4

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2674.9773
$ws.Range("I132").Value = 2628.1025
$ws.Range("J132").Value = 3040.6
$ws.Range("K132").Value = 7884.3075
$ws.Range("L132").Value = 9121.799999999999
$ws.Range("M132").Value = -5354.3075
$ws.Range("N132").Value = -14181.8
$ws.Range("H134").Value = 58590.668
$ws.Range("J134").Value = 58590.668
$ws.Range("L134").Value = 58590.668
$ws.Range("N134").Value = -68730.66800000001
$ws.Range("H136").Value = 1574.6072
$ws.Range("I136").Value = 1379.0476
$ws.Range("J136").Value = 2161.2856
$ws.Range("K136").Value = 4137.142800000001
$ws.Range("L136").Value = 6483.8568
$ws.Range("M136").Value = -1587.142800000001
$ws.Range("N136").Value = -11583.8568
$ws.Range("H137").Value = 49800
$ws.Range("J137").Value = 49800
$ws.Range("L137").Value = 49800
$ws.Range("N137").Value = -60000
$ws.Range("H141").Value = 65715
$ws.Range("J141").Value = 65715
$ws.Range("L141").Value = 65715
$ws.Range("N141").Value = -76075

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1894.375
$ws.Range("I132").Value = 1757.3
$ws.Range("K132").Value = 5271.9
$ws.Range("M132").Value = -2741.9
